# This workbook previously tracked two products (one worksheet per
# product). This change implements the "delete product" action: the
# second product's sheet ("ID_49d69c3") is removed from the workbook, and
# the remaining tracked product's row is refreshed with the latest
# scrape date.

$wb = $excel.ActiveWorkbook

# Remove the worksheet for the deleted product (ID_49d69c3).
$sheetToDelete = $wb.Worksheets.Item("ID_49d69c3")
[void]$sheetToDelete.Delete()

# Update the remaining product's date to the latest tracked date.
$ws = $wb.Worksheets.Item("ID_cd22e94")
$ws.Range("B2").Value = "17/02/2025"
